$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of row number -> final exam score (column F) to fill in.
# Rows with a "G" remark (e.g. deferred exam) or rows that already
# contain a numeric value in column F are intentionally left untouched.
$finalScores = @{
    2 = 21
    3 = 48
    4 = 35
    5 = 45
    6 = 56
    7 = 58
    8 = 59
    9 = 57
    11 = 33
    12 = 76
    13 = 41
    14 = 61
    15 = 49
    17 = 64
    18 = 63
    19 = 47
    20 = 68
    21 = 75
    22 = 52
    23 = 57
    24 = 56
    25 = 73
    26 = 67
    27 = 61
    28 = 92
    29 = 71
    30 = 95
    31 = 40
    32 = 30
    33 = 41
    34 = 46
    35 = 33
    36 = 47
    37 = 40
    38 = 45
    39 = 88
    41 = 27
    43 = 54
    44 = 65
    45 = 54
    46 = 27
    47 = 28
    48 = 41
    49 = 48
    51 = 73
    53 = 50
    55 = 70
    56 = 58
    57 = 68
    58 = 62
    59 = 58
    61 = 45
    62 = 48
    63 = 75
    65 = 45
    67 = 46
    68 = 58
    69 = 45
    70 = 33
    71 = 68
    72 = 39
    73 = 42
}

foreach ($row in $finalScores.Keys) {
    $ws.Cells.Item($row, 6).Value = $finalScores[$row]
}
